$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells removed by the naive-forecaster fix
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Update recomputed forecast values (floating point refresh)
$ws.Range("E3").Value = 4.566338461218034
$ws.Range("E4").Value = 8.262942840582976
$ws.Range("C5").Value = 8.604123301398015
$ws.Range("E5").Value = 8.260999835306748
$ws.Range("C6").Value = 6.277541464866965
$ws.Range("E6").Value = 7.915558093865016
$ws.Range("E7").Value = 4.862860110364853
$ws.Range("C8").Value = 6.535114773304795
$ws.Range("C9").Value = 6.334380382529448
$ws.Range("C11").Value = 3.889938592324382
$ws.Range("E11").Value = 3.430035192100678
$ws.Range("E12").Value = 3.624426704091555
$ws.Range("C13").Value = 2.513767348245066
$ws.Range("E14").Value = 3.771815305047843
$ws.Range("E15").Value = 2.551173534479356
$ws.Range("E18").Value = 2.56219956496937
$ws.Range("C19").Value = 2.321003614014905
$ws.Range("C20").Value = 2.508469427909921
$ws.Range("E20").Value = 2.632055757778873
$ws.Range("E23").Value = 2.874490643722782
$ws.Range("E25").Value = 1.947147525128057
$ws.Range("E26").Value = 1.467147844249128
$ws.Range("C27").Value = 1.388614840712399
$ws.Range("E28").Value = 0.9311475558545279
$ws.Range("C29").Value = 2.983246785467775
$ws.Range("C30").Value = 3.047037961814514
$ws.Range("E30").Value = 2.75705424928776
$ws.Range("E32").Value = 2.152035263856322
$ws.Range("E33").Value = 0.1082486211169131
$ws.Range("C34").Value = -0.22288476972816
$ws.Range("E34").Value = 1.7415595764392
$ws.Range("E36").Value = -0.1151106897825049
$ws.Range("C37").Value = -1.169239269654432
$ws.Range("E37").Value = -0.08704962334619148
$ws.Range("C38").Value = -1.165854108406639
$ws.Range("E39").Value = 0.611102887179138
$ws.Range("C41").Value = 2.240069601028782
$ws.Range("E42").Value = 2.055357398179125
$ws.Range("E43").Value = 1.649971980071019
$ws.Range("C45").Value = 0.8849036736221061
$ws.Range("E45").Value = 0.4552838284335392
$ws.Range("E46").Value = 0.9064937165318865
$ws.Range("C48").Value = -0.1883299148263795
$ws.Range("C49").Value = 1.978716901346367
$ws.Range("E49").Value = 1.861561193014349
$ws.Range("C50").Value = 2.039329803030099
$ws.Range("C52").Value = 3.083905204716264
$ws.Range("C53").Value = 2.3295061262957
$ws.Range("E53").Value = 1.88649514615955
